# Solicitacao de Mudanca - apply "Justificativa" paragraphs justification
# (w:jc w:val="both") and move the _GoBack bookmark so that it wraps the
# "5.2 - Justificativa" paragraph plus the blank paragraph that follows it.

$d = $word.ActiveDocument
$wdAlignParagraphJustify = 3

function Set-JustifyByFind($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $para.Format.Alignment = $wdAlignParagraphJustify
    return $para
}

# 1) "1.1 - Descricao sumaria" body paragraph
Set-JustifyByFind("se o usuário deixar campo em branco e pressiona") | Out-Null

# 2) "1.2 - Justificativa" body paragraph
Set-JustifyByFind("banco de dados terá grande erro em sua estrutura") | Out-Null

# 3) "2.1 - Descricao sumaria" body paragraph
Set-JustifyByFind("usuário poderá cadastrar idade sem nenhuma restrição") | Out-Null

# 4) "2.2 - Justificativa" body paragraph
Set-JustifyByFind("a empresa poderá pagar e sofrer um processor") | Out-Null

# 5) "3.1 - Descricao sumaria" body paragraph
Set-JustifyByFind("ao buscar o mesmo") | Out-Null

# 6) "3.2 - Justificativa" body paragraph (no pStyle, Times New Roman run)
Set-JustifyByFind("conforme a gramatica, esta escrita de forma incorreta") | Out-Null

# 7) "4.1 - Descricao sumaria" body paragraph
Set-JustifyByFind("ao excluir o mesmo") | Out-Null

# 8) "4.2 - Justificativa" blank paragraph that precedes the text paragraph
#    (no pStyle, Times New Roman run) -- find via the following text
#    paragraph (left untouched) and step backwards to the blank one.
$rng117 = $d.Content
$found117 = $rng117.Find.Execute("O usuário fica sem entender o que aconteceu", $true, $false, `
                                  $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found117) {
    throw "Find failed for paragraph 117 anchor"
}
$para117 = $rng117.Paragraphs(1)
$para116 = $para117.Previous()
$para116.Format.Alignment = $wdAlignParagraphJustify

# 9) "5.1 - Descricao sumaria" body paragraph (also has ind firstLine)
Set-JustifyByFind("onde espera um numero de") | Out-Null

# 10) "5.2 - Justificativa" body paragraph -- also where the bookmark start
#     needs to move to (beginning of the paragraph's content).
$para142 = Set-JustifyByFind("promovendo um uso adequado do")

# 11) blank paragraph right after it (before "5.3 - ..." heading) -- also
#     where the bookmark end needs to move to.
$para143 = $para142.Next()
$para143.Format.Alignment = $wdAlignParagraphJustify

# Move the _GoBack bookmark so it spans from the very start of the
# "5.2 - Justificativa" paragraph through the end of the blank paragraph
# that follows it.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$newBookmarkRange = $d.Range($para142.Range.Start, $para143.Range.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
